# Updates Coin/Link/Price/Volume(1h) cells to the latest scraped cryptos data.
# Price cells that look like plain numbers ("243.77") would otherwise be
# auto-converted to numeric values by Excel, so we force them to stay text
# (matching the source data, which is always a literal string) by switching
# the cell to Text format for the write and then restoring its original style.
function Set-TextValue {
    param($ws, $cellRef, $val)
    $prevStyle = $ws.Range($cellRef).Style
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = $prevStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22/23 swapped position: Chainlink now ranks above WrappedliquidstakedEther2.0
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"

$priceUpdates = @(
    @{ Cell = "D2"; Value = "29.797.88" },
    @{ Cell = "D3"; Value = "1.891.17" },
    @{ Cell = "D4"; Value = "0.9999" },
    @{ Cell = "D5"; Value = "0.7808" },
    @{ Cell = "D6"; Value = "243.77" },
    @{ Cell = "D7"; Value = "1.000" },
    @{ Cell = "D8"; Value = "0.3128" },
    @{ Cell = "D9"; Value = "25.28" },
    @{ Cell = "D10"; Value = "0.07167" },
    @{ Cell = "D11"; Value = "0.08051" },
    @{ Cell = "D12"; Value = "0.7624" },
    @{ Cell = "D13"; Value = "5.463" },
    @{ Cell = "D14"; Value = "1.842.90" },
    @{ Cell = "D15"; Value = "92.21" },
    @{ Cell = "D16"; Value = "6.159" },
    @{ Cell = "D17"; Value = "29.775.60" },
    @{ Cell = "D18"; Value = "13.94" },
    @{ Cell = "D19"; Value = "243.58" },
    @{ Cell = "D20"; Value = "0.000007761" },
    @{ Cell = "D21"; Value = "0.9998" },
    @{ Cell = "D22"; Value = "8.118" },
    @{ Cell = "D23"; Value = "2.119.84" },
    @{ Cell = "D24"; Value = "1.000" },
    @{ Cell = "D25"; Value = "0.1622" },
    @{ Cell = "D26"; Value = "9.388" },
    @{ Cell = "D27"; Value = "161.62" },
    @{ Cell = "D28"; Value = "18.70" },
    @{ Cell = "D29"; Value = "2.046" },
    @{ Cell = "D30"; Value = "1.413" },
    @{ Cell = "D31"; Value = "1.544" },
    @{ Cell = "D32"; Value = "4.474" },
    @{ Cell = "D33"; Value = "4.095" },
    @{ Cell = "D34"; Value = "0.05534" },
    @{ Cell = "D35"; Value = "1.263" },
    @{ Cell = "D36"; Value = "0.7436" },
    @{ Cell = "D37"; Value = "0.9956" },
    @{ Cell = "D38"; Value = "2.619" },
    @{ Cell = "D39"; Value = "0.01913" },
    @{ Cell = "D40"; Value = "2.776" },
    @{ Cell = "D41"; Value = "1.141.95" },
    @{ Cell = "D42"; Value = "73.69" },
    @{ Cell = "D43"; Value = "0.4413" },
    @{ Cell = "D44"; Value = "0.8546" },
    @{ Cell = "D45"; Value = "5.843" },
    @{ Cell = "D46"; Value = "0.9999" },
    @{ Cell = "D47"; Value = "103.57" },
    @{ Cell = "D48"; Value = "1.878" },
    @{ Cell = "D49"; Value = "9.891" },
    @{ Cell = "D50"; Value = "7.435" },
    @{ Cell = "D51"; Value = "3.016" }
)
foreach ($item in $priceUpdates) {
    Set-TextValue $ws $item.Cell $item.Value
}

$volumeUpdates = @(
    @{ Cell = "E2"; Value = "  -1.30%  " },
    @{ Cell = "E3"; Value = "  -1.12%  " },
    @{ Cell = "E4"; Value = "  -0.07%  " },
    @{ Cell = "E5"; Value = "  -4.78%  " },
    @{ Cell = "E6"; Value = "  -0.05%  " },
    @{ Cell = "E7"; Value = "  +0.00%  " },
    @{ Cell = "E8"; Value = "  -4.07%  " },
    @{ Cell = "E9"; Value = "  -6.36%  " },
    @{ Cell = "E10"; Value = "  +0.91%  " },
    @{ Cell = "E11"; Value = "  -0.37%  " },
    @{ Cell = "E12"; Value = "  -2.20%  " },
    @{ Cell = "E13"; Value = "  +2.34%  " },
    @{ Cell = "E14"; Value = "  -3.95%  " },
    @{ Cell = "E15"; Value = "  -1.79%  " },
    @{ Cell = "E16"; Value = "  +3.19%  " },
    @{ Cell = "E17"; Value = "  -1.39%  " },
    @{ Cell = "E18"; Value = "  -2.73%  " },
    @{ Cell = "E19"; Value = "  -2.11%  " },
    @{ Cell = "E20"; Value = "  -0.88%  " },
    @{ Cell = "E21"; Value = "  -0.08%  " },
    @{ Cell = "E22"; Value = "  +8.17%  " },
    @{ Cell = "E23"; Value = "  -1.76%  " },
    @{ Cell = "E24"; Value = "  -0.09%  " },
    @{ Cell = "E25"; Value = "  -3.69%  " },
    @{ Cell = "E26"; Value = "  -0.58%  " },
    @{ Cell = "E27"; Value = "  -3.89%  " },
    @{ Cell = "E28"; Value = "  -1.91%  " },
    @{ Cell = "E29"; Value = "  -3.06%  " },
    @{ Cell = "E30"; Value = "  +3.10%  " },
    @{ Cell = "E31"; Value = "  +0.84%  " },
    @{ Cell = "E32"; Value = "  +3.40%  " },
    @{ Cell = "E33"; Value = "  -0.48%  " },
    @{ Cell = "E34"; Value = "  -3.52%  " },
    @{ Cell = "E35"; Value = "  -1.15%  " },
    @{ Cell = "E36"; Value = "  +0.53%  " },
    @{ Cell = "E37"; Value = "  -0.45%  " },
    @{ Cell = "E38"; Value = "  -3.94%  " },
    @{ Cell = "E39"; Value = "  -1.11%  " },
    @{ Cell = "E40"; Value = "  -1.09%  " },
    @{ Cell = "E41"; Value = "  +11.09%  " },
    @{ Cell = "E42"; Value = "  +0.05%  " },
    @{ Cell = "E43"; Value = "  -1.38%  " },
    @{ Cell = "E44"; Value = "  +0.73%  " },
    @{ Cell = "E45"; Value = "  -1.99%  " },
    @{ Cell = "E46"; Value = "  -0.04%  " },
    @{ Cell = "E47"; Value = "  +0.65%  " },
    @{ Cell = "E48"; Value = "  -2.31%  " },
    @{ Cell = "E49"; Value = "  -0.22%  " },
    @{ Cell = "E50"; Value = "  -2.19%  " },
    @{ Cell = "E51"; Value = "  +10.90%  " }
)
foreach ($item in $volumeUpdates) {
    $ws.Range($item.Cell).Value = $item.Value
}
